$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 24 with results for Bandpower + PCA + NuSVM (linear kernel)
$ws.Range("A24").Value = "Bandpower + PCA + NuSVM (linear kernel)"
$ws.Range("B24").Value = 0.9187
$ws.Range("B24").NumberFormat = "0.00%"
$ws.Range("C24").Value = "17/19"
$ws.Range("D24").Value = "RH"
$ws.Range("E24").Value = "11, 36, 52"
$ws.Range("F24").Value = "nu=0.8, n_components=3, freq bands (Hz) 4-8,8-13,13-30, timing 375-500"

$ws.Range("F24").Select()
